# Refresh the cryptos price/volume snapshot (columns D and E) with the latest
# scraped values. D/E cells are stored as plain text (no numeric format applied
# to the sheet), so each write is apostrophe-prefixed to force a text literal -
# otherwise Excel would silently coerce number-looking strings like "19.74" into
# a numeric cell (and mangle the 2-dot "25.792.39"-style values entirely). The
# Style reset afterwards clears the "quote prefix" text format Excel applies so
# the cell keeps the workbook default style, matching the untouched cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "25.792.39"
Set-TextValue "E2" "  -0.48%  "

# Row 3
Set-TextValue "D3" "1.636.33"
Set-TextValue "E3" "  -0.36%  "

# Row 4
Set-TextValue "E4" "  -0.06%  "

# Row 5
Set-TextValue "D5" "215.66"

# Row 6
Set-TextValue "E6" "  -0.83%  "

# Row 7
Set-TextValue "E7" "  -0.08%  "

# Row 8
Set-TextValue "E8" "  -0.91%  "

# Row 9
Set-TextValue "E9" "  -1.32%  "

# Row 10
Set-TextValue "D10" "19.74"
Set-TextValue "E10" "  -2.76%  "

# Row 11
Set-TextValue "D11" "0.0791"
Set-TextValue "E11" "  +1.21%  "

# Row 12
Set-TextValue "D12" "4.28"
Set-TextValue "E12" "  +0.32%  "

# Row 13
Set-TextValue "D13" "1.863.47"
Set-TextValue "E13" "  -0.26%  "

# Row 14
Set-TextValue "D14" "1.636.35"
Set-TextValue "E14" "  -0.61%  "

# Row 15
Set-TextValue "E15" "  -0.91%  "

# Row 16
Set-TextValue "D16" "0.0₃0769"
Set-TextValue "E16" "  -0.37%  "

# Row 17
Set-TextValue "D17" "63.20"
Set-TextValue "E17" "  -0.66%  "

# Row 18
Set-TextValue "D18" "25.819.46"
Set-TextValue "E18" "  -0.39%  "

# Row 19
Set-TextValue "E19" "  -0.13%  "

# Row 20
Set-TextValue "E20" "  +1.22%  "

# Row 21
Set-TextValue "D21" "192.96"
Set-TextValue "E21" "  -1.04%  "

# Row 22
Set-TextValue "D22" "9.98"
Set-TextValue "E22" "  -0.10%  "

# Row 23
Set-TextValue "D23" "6.39"
Set-TextValue "E23" "  +2.34%  "

# Row 24
Set-TextValue "E24" "  -0.08%  "

# Row 25
Set-TextValue "D25" "1.81"
Set-TextValue "E25" "  +2.82%  "

# Row 26
Set-TextValue "E26" "  +2.62%  "

# Row 27
Set-TextValue "E27" "  -0.17%  "

# Row 28
Set-TextValue "E28" "  +1.14%  "

# Row 29
Set-TextValue "D29" "15.54"
Set-TextValue "E29" "  -0.47%  "

# Row 30
Set-TextValue "E30" "  -0.73%  "

# Row 32
Set-TextValue "D32" "3.34"
Set-TextValue "E32" "  +0.68%  "

# Row 33
Set-TextValue "E33" "  -0.64%  "

# Row 34
Set-TextValue "E34" "  -0.22%  "

# Row 35
Set-TextValue "E35" "  -0.12%  "

# Row 36
Set-TextValue "D36" "0.906"
Set-TextValue "E36" "  -0.54%  "

# Row 37
Set-TextValue "D37" "1.134.67"
Set-TextValue "E37" "  +0.43%  "

# Row 38
Set-TextValue "E38" "  -2.07%  "

# Row 39
Set-TextValue "D39" "0.545"
Set-TextValue "E39" "  -1.54%  "

# Row 40
Set-TextValue "E40" "  -0.96%  "

# Row 41
Set-TextValue "E41" "  +0.64%  "

# Row 42
Set-TextValue "D42" "5.56"
Set-TextValue "E42" "  +1.07%  "

# Row 43
Set-TextValue "D43" "100.35"
Set-TextValue "E43" "  +0.28%  "

# Row 44
Set-TextValue "D44" "0.806"
Set-TextValue "E44" "  +0.47%  "

# Row 45
Set-TextValue "D45" "1.773.08"
Set-TextValue "E45" "  -0.40%  "

# Row 46
Set-TextValue "D46" "0.0₆0110"
Set-TextValue "E46" "  -0.40%  "

# Row 47
Set-TextValue "D47" "55.29"
Set-TextValue "E47" "  -0.94%  "

# Row 48
Set-TextValue "E48" "  -1.51%  "

# Row 49
Set-TextValue "E49" "  -0.25%  "

# Row 51
Set-TextValue "D51" "7.51"
Set-TextValue "E51" "  -2.86%  "
